$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.432.39"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "3.696.17"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'693.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'162.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "3.694.49"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -4.65%  "
$ws.Range("E10").Value = "  -8.10%  "
$ws.Range("D11").Value = "'7.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("E13").Value = "  -5.18%  "
$ws.Range("D14").Value = "'33.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.26%  "
$ws.Range("D15").Value = "4.315.74"
$ws.Range("E15").Value = "  -3.10%  "
$ws.Range("D16").Value = "3.695.15"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "69.462.47"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'16.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.22%  "
$ws.Range("D20").Value = "'6.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.65%  "
$ws.Range("D21").Value = "'481.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.23%  "
$ws.Range("D22").Value = "'9.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.28%  "
$ws.Range("E23").Value = "  -7.31%  "
$ws.Range("D24").Value = "'80.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("D25").Value = "3.840.25"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  -9.57%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -5.97%  "
$ws.Range("D29").Value = "'9.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.81%  "
$ws.Range("E30").Value = "  -10.82%  "
$ws.Range("E31").Value = "  -10.25%  "
$ws.Range("D32").Value = "'6.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.05%  "
$ws.Range("E33").Value = "  -7.81%  "
$ws.Range("D34").Value = "'27.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.18%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "3.663.28"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").Value = "'8.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.63%  "
$ws.Range("D39").Value = "'6.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").Value = "'2.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").Value = "'0.0931"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.87%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'0.955"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.36%  "
$ws.Range("D45").Value = "'164.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("D46").Value = "'48.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "'30.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("E48").Value = "  -15.46%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "'0.000284"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.75%  "
